# Add two new columns (I: "I0", J: "IF") to the header row, copying the
# header formatting (bold font, border, centered alignment) from the
# existing "IP" header cell (H1), then fill in the two data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy H1's formatting onto I1:J1 (xlPasteFormats) before writing values,
# so the new header cells pick up the same style used by the other headers.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 5

$ws.Range("I3").Value = 3
$ws.Range("J3").Value = 3
